# fix scaling of 2024 logs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 25
$ws.Range("B2").Value = 57
$ws.Range("B3").Value = 128
$ws.Range("B4").Value = 247
